# Fruta / hortaliza, semanal
#
# Inserts a new weekly price record (3 rows: Especial/Primera/Segunda) at the
# top of the "Terminal La Palmera de La Serena - Mango" history block
# (previously starting at row 988), pushing the existing 39 historical rows
# (988-1026) down by 3 rows (to 991-1029).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 988 (shifts old rows 988.. down to 991..)
$ws.Range("A988:A990").EntireRow.Insert()

$qualities = @("Especial", "Primera", "Segunda")

for ($i = 0; $i -lt 3; $i++) {
    $r = 988 + $i

    $ws.Cells.Item($r, 1).Value  = 8
    $ws.Cells.Item($r, 2).Value  = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value  = "Coquimbo"
    $ws.Cells.Item($r, 4).Value  = 44939
    $ws.Cells.Item($r, 5).Value  = 4
    $ws.Cells.Item($r, 6).Value  = "Fruta"
    $ws.Cells.Item($r, 7).Value  = 100108
    $ws.Cells.Item($r, 8).Value  = "Tropicales y subtropicales"
    $ws.Cells.Item($r, 9).Value  = 100108002
    $ws.Cells.Item($r, 10).Value = "Mango"
    $ws.Cells.Item($r, 11).Value = "Sin especificar"
    $ws.Cells.Item($r, 12).Value = $qualities[$i]
    $ws.Cells.Item($r, 13).Value = 512
    $ws.Cells.Item($r, 14).Value = 6000
    $ws.Cells.Item($r, 15).Value = 6500
    $ws.Cells.Item($r, 16).Value = 6250
    $ws.Cells.Item($r, 17).Value = "$/bandeja 4 kilos"
    $ws.Cells.Item($r, 18).Value = "Perú"
    $ws.Cells.Item($r, 19).Value = 1562
    $ws.Cells.Item($r, 20).Value = 4
}

Write-Output "Inserted rows 988-990; dimension now A1:T1029"
